$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "INTRON/EXON" front-end-task note (row 8, column E) with the
# expanded explanation that now also includes the Perl code snippet for
# the front layer.
$ws.Range("E8").Value = 'This task can be completed by only giving the front end the full DNA sequence and an hash of exons positions/length; both retrievable with queries; all the front end needs is to know which sbstring within the main string to highligh; same we decided to do for the restriction sites task. Code for the front layer:                                                                 foreach my $key (keys %exons)   {
   substr($sequence, $key, $exons{$key}) = "<div>" . substr($sequence,    $key, $exons{$key} ) . "</div>";
}
print "<p> $sequence </p>";'

# The much longer note no longer fits the old row height, so the row grows
# to accommodate the wrapped text.
$ws.Rows.Item(8).RowHeight = 165

# Move the current selection down to the newly edited row.
[void]$ws.Range("F8").Select()
